$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("I2").Value = 0.07105816984606554
$ws.Range("J2").Value = 0.07105816984606554
$ws.Range("M2").Value = 68.637375
$ws.Range("N2").Value = 205.912125
$ws.Range("O2").Value = 0.5415701538216162
$ws.Range("P2").Value = 0.5415701538216162
$ws.Range("Q2").Value = 2.447265605625
$ws.Range("R2").Value = 22.025390450625
$ws.Range("S2").Value = 0.03848298397381625
$ws.Range("T2").Value = 0.03848298397381625

# Row 3
$ws.Range("I3").Value = 0.07105816984606554
$ws.Range("J3").Value = 0.07105816984606554
$ws.Range("O3").Value = 0.08718851262838957
$ws.Range("P3").Value = 0.08718851262838957
$ws.Range("S3").Value = 0.006195456138973937
$ws.Range("T3").Value = 0.006195456138973937

# Row 4
$ws.Range("I4").Value = 0.07105816984606554
$ws.Range("J4").Value = 0.07105816984606554
$ws.Range("M4").Value = 16.21089566666667
$ws.Range("N4").Value = 48.632687
$ws.Range("O4").Value = 0.1279089892319285
$ws.Range("P4").Value = 0.1279089892319285
$ws.Range("Q4").Value = 0.5779994849950001
$ws.Range("R4").Value = 5.201995364955001
$ws.Range("S4").Value = 0.009088978681680944
$ws.Range("T4").Value = 0.009088978681680944

# Row 5
$ws.Range("I5").Value = 0.07105816984606554
$ws.Range("J5").Value = 0.07105816984606554
$ws.Range("M5").Value = 20.32546233333333
$ws.Range("N5").Value = 60.976387
$ws.Range("O5").Value = 0.1603741949973873
$ws.Range("P5").Value = 0.1603741949973873
$ws.Range("Q5").Value = 0.724704359495
$ws.Range("R5").Value = 6.522339235455001
$ws.Range("S5").Value = 0.01139589678705038
$ws.Range("T5").Value = 0.01139589678705038

# Row 6
$ws.Range("I6").Value = 0.07105816984606554
$ws.Range("J6").Value = 0.07105816984606554
$ws.Range("M6").Value = 10.513928
$ws.Range("N6").Value = 31.541784
$ws.Range("O6").Value = 0.08295814932067838
$ws.Range("P6").Value = 0.08295814932067838
$ws.Range("Q6").Value = 0.37487410284
$ws.Range("R6").Value = 3.37386692556
$ws.Range("S6").Value = 0.005894854264544031
$ws.Range("T6").Value = 0.005894854264544031

# Row 7
$ws.Range("G7").Value = 0.1910563333333334
$ws.Range("H7").Value = 0.573169
$ws.Range("I7").Value = 0.3807632417379475
$ws.Range("J7").Value = 0.3807632417379474
$ws.Range("M7").Value = 68.637375
$ws.Range("N7").Value = 205.912125
$ws.Range("O7").Value = 0.5415701538216162
$ws.Range("P7").Value = 0.5415701538216162
$ws.Range("Q7").Value = 13.113605197125
$ws.Range("R7").Value = 118.022446774125
$ws.Range("S7").Value = 0.2062100073976375
$ws.Range("T7").Value = 0.2062100073976374

# Row 8
$ws.Range("G8").Value = 0.1910563333333334
$ws.Range("H8").Value = 0.573169
$ws.Range("I8").Value = 0.3807632417379475
$ws.Range("J8").Value = 0.3807632417379474
$ws.Range("O8").Value = 0.08718851262838957
$ws.Range("P8").Value = 0.08718851262838957
$ws.Range("S8").Value = 0.03319818071069559
$ws.Range("T8").Value = 0.03319818071069558

# Row 9
$ws.Range("G9").Value = 0.1910563333333334
$ws.Range("H9").Value = 0.573169
$ws.Range("I9").Value = 0.3807632417379475
$ws.Range("J9").Value = 0.3807632417379474
$ws.Range("M9").Value = 16.21089566666667
$ws.Range("N9").Value = 48.632687
$ws.Range("O9").Value = 0.1279089892319285
$ws.Range("P9").Value = 0.1279089892319285
$ws.Range("Q9").Value = 3.097194286122556
$ws.Range("R9").Value = 27.87474857510301
$ws.Range("S9").Value = 0.04870304138737331
$ws.Range("T9").Value = 0.0487030413873733

# Row 10
$ws.Range("G10").Value = 0.1910563333333334
$ws.Range("H10").Value = 0.573169
$ws.Range("I10").Value = 0.3807632417379475
$ws.Range("J10").Value = 0.3807632417379474
$ws.Range("M10").Value = 20.32546233333333
$ws.Range("N10").Value = 60.976387
$ws.Range("O10").Value = 0.1603741949973873
$ws.Range("P10").Value = 0.1603741949973873
$ws.Range("Q10").Value = 3.883308306711445
$ws.Range("R10").Value = 34.949774760403
$ws.Range("S10").Value = 0.06106459837831892
$ws.Range("T10").Value = 0.06106459837831891

# Row 11
$ws.Range("G11").Value = 0.1910563333333334
$ws.Range("H11").Value = 0.573169
$ws.Range("I11").Value = 0.3807632417379475
$ws.Range("J11").Value = 0.3807632417379474
$ws.Range("M11").Value = 10.513928
$ws.Range("N11").Value = 31.541784
$ws.Range("O11").Value = 0.08295814932067838
$ws.Range("P11").Value = 0.08295814932067838
$ws.Range("Q11").Value = 2.008752532610667
$ws.Range("R11").Value = 18.078772793496
$ws.Range("S11").Value = 0.03158741386392221
$ws.Range("T11").Value = 0.0315874138639222

# Row 12
$ws.Range("G12").Value = 0.2750606666666667
$ws.Range("H12").Value = 0.825182
$ws.Range("I12").Value = 0.5481785884159871
$ws.Range("J12").Value = 0.548178588415987
$ws.Range("M12").Value = 68.637375
$ws.Range("N12").Value = 205.912125
$ws.Range("O12").Value = 0.5415701538216162
$ws.Range("P12").Value = 0.5415701538216162
$ws.Range("Q12").Value = 18.87944212575
$ws.Range("R12").Value = 169.91497913175
$ws.Range("S12").Value = 0.2968771624501626
$ws.Range("T12").Value = 0.2968771624501625

# Row 13
$ws.Range("G13").Value = 0.2750606666666667
$ws.Range("H13").Value = 0.825182
$ws.Range("I13").Value = 0.5481785884159871
$ws.Range("J13").Value = 0.548178588415987
$ws.Range("O13").Value = 0.08718851262838957
$ws.Range("P13").Value = 0.08718851262838957
$ws.Range("Q13").Value = 3.039440904529778
$ws.Range("R13").Value = 27.354968140768
$ws.Range("S13").Value = 0.04779487577872006
$ws.Range("T13").Value = 0.04779487577872005

# Row 14
$ws.Range("G14").Value = 0.2750606666666667
$ws.Range("H14").Value = 0.825182
$ws.Range("I14").Value = 0.5481785884159871
$ws.Range("J14").Value = 0.548178588415987
$ws.Range("M14").Value = 16.21089566666667
$ws.Range("N14").Value = 48.632687
$ws.Range("O14").Value = 0.1279089892319285
$ws.Range("P14").Value = 0.1279089892319285
$ws.Range("Q14").Value = 4.458979769337112
$ws.Range("R14").Value = 40.13081792403401
$ws.Range("S14").Value = 0.07011696916287426
$ws.Range("T14").Value = 0.07011696916287424

# Row 15
$ws.Range("G15").Value = 0.2750606666666667
$ws.Range("H15").Value = 0.825182
$ws.Range("I15").Value = 0.5481785884159871
$ws.Range("J15").Value = 0.548178588415987
$ws.Range("M15").Value = 20.32546233333333
$ws.Range("N15").Value = 60.976387
$ws.Range("O15").Value = 0.1603741949973873
$ws.Range("P15").Value = 0.1603741949973873
$ws.Range("Q15").Value = 5.590735219714889
$ws.Range("R15").Value = 50.316616977434
$ws.Range("S15").Value = 0.08791369983201805
$ws.Range("T15").Value = 0.08791369983201804

# Row 16
$ws.Range("G16").Value = 0.2750606666666667
$ws.Range("H16").Value = 0.825182
$ws.Range("I16").Value = 0.5481785884159871
$ws.Range("J16").Value = 0.548178588415987
$ws.Range("M16").Value = 10.513928
$ws.Range("N16").Value = 31.541784
$ws.Range("O16").Value = 0.08295814932067838
$ws.Range("P16").Value = 0.08295814932067838
$ws.Range("Q16").Value = 2.891968044965334
$ws.Range("R16").Value = 26.027712404688
$ws.Range("S16").Value = 0.04547588119221215
$ws.Range("T16").Value = 0.04547588119221214
